$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("operation class")

# Row 2 - new_sssi_e
$ws.Range("B2").Value = "6 (1.52%)"
$ws.Range("C2").Value = "2 (0.97%)"
$ws.Range("E2").Value = "0 (0.00%)"
$ws.Range("F2").Value = "1 (0.71%)"
$ws.Range("G2").Value = "0 (0.00%)"

# Row 3 - dehis_e
$ws.Range("E3").Value = "0 (0.00%)"
$ws.Range("F3").Value = "14 (10.00%)"
$ws.Range("G3").Value = "0 (0.00%)"

# Row 6 - urninfec_e
$ws.Range("B6").Value = "1 (0.25%)"
$ws.Range("C6").Value = "6 (2.91%)"
$ws.Range("D6").Value = "1 (3.23%)"
$ws.Range("E6").Value = "0 (0.00%)"
$ws.Range("F6").Value = "4 (2.86%)"
$ws.Range("G6").Value = "0 (0.00%)"

# Row 10 - any_complication
$ws.Range("B10").Value = "7 (1.77%)"
$ws.Range("C10").Value = "7 (3.40%)"
$ws.Range("D10").Value = "1 (3.23%)"
$ws.Range("E10").Value = "0 (0.00%)"
$ws.Range("F10").Value = "19 (13.57%)"
$ws.Range("G10").Value = "0 (0.00%)"
